$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 8.371231
$ws.Range("H2").Value = 25.113693
$ws.Range("I2").Value = 0.1018690981095697
$ws.Range("J2").Value = 0.1018690981095697
$ws.Range("M2").Value = 2.769264333333334
$ws.Range("N2").Value = 8.307793
$ws.Range("O2").Value = 0.1388016358751757
$ws.Range("P2").Value = 0.1388016358751757
$ws.Range("Q2").Value = 23.18215143439433
$ws.Range("R2").Value = 208.639362909549
$ws.Range("S2").Value = 0.01413959746273705
$ws.Range("T2").Value = 0.01413959746273705
$ws.Range("G3").Value = 8.371231
$ws.Range("H3").Value = 25.113693
$ws.Range("I3").Value = 0.1018690981095697
$ws.Range("J3").Value = 0.1018690981095697
$ws.Range("M3").Value = 1.484487666666667
$ws.Range("N3").Value = 4.453463
$ws.Range("O3").Value = 0.07440579582442265
$ws.Range("P3").Value = 0.07440579582442265
$ws.Range("Q3").Value = 12.42698917431767
$ws.Range("R3").Value = 111.842902568859
$ws.Range("S3").Value = 0.007579651314758723
$ws.Range("T3").Value = 0.007579651314758723
$ws.Range("G4").Value = 8.371231
$ws.Range("H4").Value = 25.113693
$ws.Range("I4").Value = 0.1018690981095697
$ws.Range("J4").Value = 0.1018690981095697
$ws.Range("M4").Value = 0.1509683333333333
$ws.Range("N4").Value = 0.452905
$ws.Range("O4").Value = 0.007566865820567083
$ws.Range("P4").Value = 0.007566865820567083
$ws.Range("Q4").Value = 1.263790792018333
$ws.Range("R4").Value = 11.374117128165
$ws.Range("S4").Value = 0.0007708297966572979
$ws.Range("T4").Value = 0.0007708297966572979
$ws.Range("G5").Value = 8.371231
$ws.Range("H5").Value = 25.113693
$ws.Range("I5").Value = 0.1018690981095697
$ws.Range("J5").Value = 0.1018690981095697
$ws.Range("M5").Value = 15.54651666666667
$ws.Range("N5").Value = 46.63955
$ws.Range("O5").Value = 0.7792257024798346
$ws.Range("P5").Value = 0.7792257024798346
$ws.Range("Q5").Value = 130.1434822620167
$ws.Range("R5").Value = 1171.29134035815
$ws.Range("S5").Value = 0.07937901953541664
$ws.Range("T5").Value = 0.07937901953541664
$ws.Range("I6").Value = 0.683327746432814
$ws.Range("J6").Value = 0.683327746432814
$ws.Range("M6").Value = 2.769264333333334
$ws.Range("N6").Value = 8.307793
$ws.Range("O6").Value = 0.1388016358751757
$ws.Range("P6").Value = 0.1388016358751757
$ws.Range("Q6").Value = 155.5035588917301
$ws.Range("R6").Value = 1399.532030025571
$ws.Range("S6").Value = 0.09484700904377187
$ws.Range("T6").Value = 0.09484700904377187
$ws.Range("I7").Value = 0.683327746432814
$ws.Range("J7").Value = 0.683327746432814
$ws.Range("M7").Value = 1.484487666666667
$ws.Range("N7").Value = 4.453463
$ws.Range("O7").Value = 0.07440579582442265
$ws.Range("P7").Value = 0.07440579582442265
$ws.Range("Q7").Value = 83.35900351545122
$ws.Range("R7").Value = 750.231031639061
$ws.Range("S7").Value = 0.05084354478224282
$ws.Range("T7").Value = 0.05084354478224282
$ws.Range("I8").Value = 0.683327746432814
$ws.Range("J8").Value = 0.683327746432814
$ws.Range("M8").Value = 0.1509683333333333
$ws.Range("N8").Value = 0.452905
$ws.Range("O8").Value = 0.007566865820567083
$ws.Range("P8").Value = 0.007566865820567083
$ws.Range("Q8").Value = 8.477382541892778
$ws.Range("R8").Value = 76.29644287703501
$ws.Range("S8").Value = 0.005170649368727591
$ws.Range("T8").Value = 0.005170649368727591
$ws.Range("I9").Value = 0.683327746432814
$ws.Range("J9").Value = 0.683327746432814
$ws.Range("M9").Value = 15.54651666666667
$ws.Range("N9").Value = 46.63955
$ws.Range("O9").Value = 0.7792257024798346
$ws.Range("P9").Value = 0.7792257024798346
$ws.Range("Q9").Value = 872.9894943348722
$ws.Range("R9").Value = 7856.905449013851
$ws.Range("S9").Value = 0.5324665432380717
$ws.Range("T9").Value = 0.5324665432380717
$ws.Range("G10").Value = 16.77784
$ws.Range("H10").Value = 50.33351999999999
$ws.Range("I10").Value = 0.2041687093602677
$ws.Range("J10").Value = 0.2041687093602677
$ws.Range("M10").Value = 2.769264333333334
$ws.Range("N10").Value = 8.307793
$ws.Range("O10").Value = 0.1388016358751757
$ws.Range("P10").Value = 0.1388016358751757
$ws.Range("Q10").Value = 46.46227390237333
$ws.Range("R10").Value = 418.1604651213599
$ws.Range("S10").Value = 0.02833895085372846
$ws.Range("T10").Value = 0.02833895085372846
$ws.Range("G11").Value = 16.77784
$ws.Range("H11").Value = 50.33351999999999
$ws.Range("I11").Value = 0.2041687093602677
$ws.Range("J11").Value = 0.2041687093602677
$ws.Range("M11").Value = 1.484487666666667
$ws.Range("N11").Value = 4.453463
$ws.Range("O11").Value = 0.07440579582442265
$ws.Range("P11").Value = 0.07440579582442265
$ws.Range("Q11").Value = 24.90649655330666
$ws.Range("R11").Value = 224.15846897976
$ws.Range("S11").Value = 0.01519133530239597
$ws.Range("T11").Value = 0.01519133530239597
$ws.Range("G12").Value = 16.77784
$ws.Range("H12").Value = 50.33351999999999
$ws.Range("I12").Value = 0.2041687093602677
$ws.Range("J12").Value = 0.2041687093602677
$ws.Range("M12").Value = 0.1509683333333333
$ws.Range("N12").Value = 0.452905
$ws.Range("O12").Value = 0.007566865820567083
$ws.Range("P12").Value = 0.007566865820567083
$ws.Range("Q12").Value = 2.532922541733333
$ws.Range("R12").Value = 22.7963028756
$ws.Range("S12").Value = 0.001544917228487504
$ws.Range("T12").Value = 0.001544917228487504
$ws.Range("G13").Value = 16.77784
$ws.Range("H13").Value = 50.33351999999999
$ws.Range("I13").Value = 0.2041687093602677
$ws.Range("J13").Value = 0.2041687093602677
$ws.Range("M13").Value = 15.54651666666667
$ws.Range("N13").Value = 46.63955
$ws.Range("O13").Value = 0.7792257024798346
$ws.Range("P13").Value = 0.7792257024798346
$ws.Range("Q13").Value = 260.8369691906667
$ws.Range("R13").Value = 2347.532722716
$ws.Range("S13").Value = 0.1590935059756558
$ws.Range("T13").Value = 0.1590935059756558
$ws.Range("G14").Value = 0.8738999999999999
$ws.Range("H14").Value = 2.6217
$ws.Range("I14").Value = 0.01063444609734852
$ws.Range("J14").Value = 0.01063444609734852
$ws.Range("M14").Value = 2.769264333333334
$ws.Range("N14").Value = 8.307793
$ws.Range("O14").Value = 0.1388016358751757
$ws.Range("P14").Value = 0.1388016358751757
$ws.Range("Q14").Value = 2.4200601009
$ws.Range("R14").Value = 21.7805409081
$ws.Range("S14").Value = 0.001476078514938353
$ws.Range("T14").Value = 0.001476078514938353
$ws.Range("G15").Value = 0.8738999999999999
$ws.Range("H15").Value = 2.6217
$ws.Range("I15").Value = 0.01063444609734852
$ws.Range("J15").Value = 0.01063444609734852
$ws.Range("M15").Value = 1.484487666666667
$ws.Range("N15").Value = 4.453463
$ws.Range("O15").Value = 0.07440579582442265
$ws.Range("P15").Value = 0.07440579582442265
$ws.Range("Q15").Value = 1.2972937719
$ws.Range("R15").Value = 11.6756439471
$ws.Range("S15").Value = 0.0007912644250251424
$ws.Range("T15").Value = 0.0007912644250251424
$ws.Range("G16").Value = 0.8738999999999999
$ws.Range("H16").Value = 2.6217
$ws.Range("I16").Value = 0.01063444609734852
$ws.Range("J16").Value = 0.01063444609734852
$ws.Range("M16").Value = 0.1509683333333333
$ws.Range("N16").Value = 0.452905
$ws.Range("O16").Value = 0.007566865820567083
$ws.Range("P16").Value = 0.007566865820567083
$ws.Range("Q16").Value = 0.1319312265
$ws.Range("R16").Value = 1.1873810385
$ws.Range("S16").Value = 0.00008046942669468954
$ws.Range("T16").Value = 0.00008046942669468954
$ws.Range("G17").Value = 0.8738999999999999
$ws.Range("H17").Value = 2.6217
$ws.Range("I17").Value = 0.01063444609734852
$ws.Range("J17").Value = 0.01063444609734852
$ws.Range("M17").Value = 15.54651666666667
$ws.Range("N17").Value = 46.63955
$ws.Range("O17").Value = 0.7792257024798346
$ws.Range("P17").Value = 0.7792257024798346
$ws.Range("Q17").Value = 13.586100915
$ws.Range("R17").Value = 122.274908235
$ws.Range("S17").Value = 0.008286633730690338
$ws.Range("T17").Value = 0.008286633730690338
